$d = $word.ActiveDocument

# --- 1. Merge runs that were split mid-word so the final text reads as single runs ---
$d.Content.Find.Execute(
    "following the instructions.  ", $true, $false, $false, $false, $false,
    $true, 1, $false, "following the instructions.  ", 2) | Out-Null

$d.Content.Find.Execute(
    "itself and selecting that option);", $true, $false, $false, $false, $false,
    $true, 1, $false, "itself and selecting that option);", 2) | Out-Null

$d.Content.Find.Execute(
    "files pane, e.g. ‘1-intro.R’", $true, $false, $false, $false, $false,
    $true, 1, $false, "files pane, e.g. ‘1-intro.R’", 2) | Out-Null

# --- 2. Insert a new introductory paragraph before the current first paragraph ---
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphBefore()

# First sentence part - bold
$r1 = $d.Range(0, 0)
$r1.Text = "The following instructions are for software installation "
$r1.Font.Bold = 1
$r1.Font.BoldBi = 1

# "after " - bold and italic
$pos = $r1.End
$r2 = $d.Range($pos, $pos)
$r2.InsertAfter("after ")
$r2.Font.Bold = 1
$r2.Font.BoldBi = 1
$r2.Font.Italic = 1

# Remainder - bold, italic explicitly turned back off
$pos2 = $r2.End
$r3 = $d.Range($pos2, $pos2)
$r3.InsertAfter("the course, if you would like to continue learning. You are not required to download or install anything before the course.")
$r3.Font.Bold = 1
$r3.Font.BoldBi = 1
$r3.Font.Italic = 0

Write-Output "Edit complete"
